# Insert a new numbered list item "código identificador de editorial"
# right after the existing "id detalle de compra" list item, matching the
# same paragraph/run formatting (numId=1 bulleted list, sz=24, etc.).

$d = $word.ActiveDocument

$anchorText = "id detalle de compra"
$newText = "código identificador de editorial"

$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq $anchorText) {
        $anchorPara = $p
        break
    }
}

if ($anchorPara -ne $null) {
    # Insert a brand-new paragraph right after the anchor; Word clones the
    # anchor's paragraph + run formatting (numbering, spacing, indents, etc.)
    # for the freshly inserted paragraph.
    $anchorPara.Range.InsertParagraphAfter()

    $newPara = $anchorPara.Next()
    $newPara.Range.Text = $newText
}
